# Regenerate the "K" column (strikeouts) values in Sheet1, column G,
# rows 2-44, replacing the previous "Strike#"-derived values with the
# newly computed K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 5
    3  = 2
    4  = 0
    5  = 4
    6  = 5
    7  = 0
    8  = 1
    9  = 1
    10 = 2
    11 = 3
    12 = 0
    13 = 3
    14 = 3
    15 = 0
    16 = 2
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 1
    23 = 4
    24 = 0
    25 = 0
    26 = 0
    27 = 1
    28 = 0
    29 = 0
    30 = 1
    31 = 0
    32 = 2
    33 = 1
    34 = 1
    35 = 1
    36 = 0
    37 = 0
    38 = 0
    39 = 1
    40 = 3
    41 = 1
    42 = 2
    43 = 0
    44 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
